$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 127, shifting rows 127:191 down to 128:192
$ws.Rows.Item(127).Insert()

# Match the date-style formatting used by column D in the rest of the table
$ws.Cells.Item(127, 4).NumberFormat = $ws.Cells.Item(128, 4).NumberFormat

# Populate the newly inserted row 127 with its data
$ws.Cells.Item(127, 1).Value = 10
$ws.Cells.Item(127, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(127, 3).Value = "La Araucanía"
$ws.Cells.Item(127, 4).Value = 44813
$ws.Cells.Item(127, 5).Value = 9
$ws.Cells.Item(127, 6).Value = 100114007
$ws.Cells.Item(127, 7).Value = "Jengibre"
$ws.Cells.Item(127, 8).Value = "Sin especificar"
$ws.Cells.Item(127, 9).Value = "Primera"
$ws.Cells.Item(127, 10).Value = 40
$ws.Cells.Item(127, 11).Value = 20000
$ws.Cells.Item(127, 12).Value = 20000
$ws.Cells.Item(127, 13).Value = 20000
$ws.Cells.Item(127, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(127, 15).Value = "Perú"
$ws.Cells.Item(127, 16).Value = 1538
$ws.Cells.Item(127, 17).Value = 13
$ws.Cells.Item(127, 18).Value = "Hortaliza"
